# The commit adds a new weekly price record for "Puerro" (Vega Modelo de
# Temuco) on 2022-11-10 (serial 44875). In the sheet's data the rows are
# kept in reverse-chronological order, so the new record is inserted as
# row 188, pushing the former rows 188-250 down to 189-251.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at 188; Excel shifts rows 188:250 down to 189:251
# and extends the used range to A1:R251 automatically.
$ws.Rows.Item(188).Insert()

# Populate the new row with the new observation. All of the
# "descriptive" columns (market, region, product, variety, quality,
# unit, origin, classification...) match the rest of the block.
$ws.Range("A188").Value = 10
$ws.Range("B188").Value = "Vega Modelo de Temuco"
$ws.Range("C188").Value = "La Araucanía"
$ws.Range("D188").Value = 44875
$ws.Range("E188").Value = 9
$ws.Range("F188").Value = 100112005
$ws.Range("G188").Value = "Puerro"
$ws.Range("H188").Value = "Azul de Maquehue"
$ws.Range("I188").Value = "Primera"
$ws.Range("J188").Value = 95
$ws.Range("K188").Value = 18000
$ws.Range("L188").Value = 18000
$ws.Range("M188").Value = 18000
$ws.Range("N188").Value = "`$/docena de paquetes"
$ws.Range("O188").Value = "Provincia de Cautín"
$ws.Range("P188").Value = 1500
$ws.Range("Q188").Value = 12
$ws.Range("R188").Value = "Hortaliza"

# Keep the date column formatted the same way as the rest of column D.
$ws.Range("D188").NumberFormat = $ws.Range("D189").NumberFormat
